# Applies the cryptos list price/volume refresh described in the commit
# "Updated cryptos list on Fri Jun 21 13:31:42 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.612.57'
$ws.Range("E2").Value = '  -2.67%  '
$ws.Range("D3").Value = '3.476.16'
$ws.Range("E3").Value = '  -1.88%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D5").Value = '580.12'
$ws.Range("E5").Value = '  -3.08%  '
$ws.Range("D6").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D6").Value = '130.09'
$ws.Range("E6").Value = '  -4.02%  '
$ws.Range("D7").Value = '3.476.75'
$ws.Range("E7").Value = '  -1.69%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D9").Value = '0.487'
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("D10").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D10").Value = '0.122'
$ws.Range("E10").Value = '  -1.09%  '
$ws.Range("D11").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D11").Value = '7.16'
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D12").Value = '0.379'
$ws.Range("E12").Value = '  -2.15%  '
$ws.Range("D13").Value = '4.083.39'
$ws.Range("E13").Value = '  -2.09%  '
$ws.Range("D14").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D14").Value = '27.26'
$ws.Range("E14").Value = '  -1.04%  '
$ws.Range("E15").Value = '  +1.33%  '
$ws.Range("D16").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D16").Value = '0.0000176'
$ws.Range("E16").Value = '  -3.18%  '
$ws.Range("D17").Value = '3.458.42'
$ws.Range("E17").Value = '  -2.68%  '
$ws.Range("D18").Value = '63.782.67'
$ws.Range("E18").Value = '  -2.58%  '
$ws.Range("D19").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D19").Value = '10.01'
$ws.Range("D20").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D20").Value = '14.30'
$ws.Range("E20").Value = '  -1.20%  '
$ws.Range("D21").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D21").Value = '5.62'
$ws.Range("E21").Value = '  -2.14%  '
$ws.Range("D22").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D22").Value = '382.06'
$ws.Range("E22").Value = '  -2.67%  '
$ws.Range("D23").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D23").Value = '0.573'
$ws.Range("E23").Value = '  -1.05%  '
$ws.Range("D24").Value = '3.618.40'
$ws.Range("E24").Value = '  -2.14%  '
$ws.Range("D25").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D25").Value = '72.93'
$ws.Range("E25").Value = '  -2.53%  '
$ws.Range("E26").Value = '  +0.26%  '
$ws.Range("D27").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D27").Value = '0.0000110'
$ws.Range("E27").Value = '  -4.13%  '
$ws.Range("D28").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D28").Value = '1.55'
$ws.Range("E28").Value = '  -4.52%  '
$ws.Range("D29").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D30").Value = '7.35'
$ws.Range("E30").Value = '  -6.81%  '
$ws.Range("D31").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D31").Value = '2.22'
$ws.Range("E31").Value = '  -3.44%  '
$ws.Range("D32").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D32").Value = '8.14'
$ws.Range("E32").Value = '  -4.26%  '
$ws.Range("D33").Value = '3.482.50'
$ws.Range("E33").Value = '  -2.08%  '
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("D35").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D35").Value = '23.54'
$ws.Range("E35").Value = '  -2.65%  '
$ws.Range("D36").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D36").Value = '0.142'
$ws.Range("E36").Value = '  -2.82%  '
$ws.Range("D37").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D37").Value = '5.19'
$ws.Range("E37").Value = '  -1.80%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D38").Value = '167.49'
$ws.Range("E38").Value = '  -0.91%  '
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D39").Value = '6.83'
$ws.Range("E39").Value = '  -0.54%  '
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D40").Value = '1.54'
$ws.Range("E40").Value = '  -2.76%  '
$ws.Range("D41").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D41").Value = '0.0793'
$ws.Range("E41").Value = '  -4.00%  '
$ws.Range("D42").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D42").Value = '26.60'
$ws.Range("E42").Value = '  +1.81%  '
$ws.Range("E43").Value = '  -2.18%  '
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("D45").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D45").Value = '41.39'
$ws.Range("E46").Value = '  -4.26%  '
$ws.Range("D47").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D47").Value = '4.35'
$ws.Range("E47").Value = '  -2.35%  '
$ws.Range("D48").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D48").Value = '1.62'
$ws.Range("E48").Value = '  -3.22%  '
$ws.Range("D49").Value = '2.427.64'
$ws.Range("E49").Value = '  +0.76%  '
$ws.Range("D50").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D50").Value = '6.80'
$ws.Range("E50").Value = '  -1.08%  '
$ws.Range("D51").NumberFormat = "@"   # keep numeric-looking text as a string
$ws.Range("D51").Value = '0.882'
$ws.Range("E51").Value = '  -1.32%  '
